$wb = $excel.ActiveWorkbook

# Sheets: 1=Personnes, 2=Articles, 3=Relations, 4=Debunks
$wsArticles = $wb.Worksheets.Item(2)
$wsRelations = $wb.Worksheets.Item(3)
$wsDebunks = $wb.Worksheets.Item(4)

# --- Articles sheet: add new row 10 for the Marie-Estelle Dupont article ---
$wsArticles.Range("A10").Value = 9
$wsArticles.Range("B10").Value = "Marie-Estelle Dupont : psycho couacs à droite toute"
$wsArticles.Range("E10").Value = "https://www.arretsurimages.net/articles/marie-estelle-dupont-psycho-couacs-a-droite-toute"
$wsArticles.Hyperlinks.Add($wsArticles.Range("E10"), "https://www.arretsurimages.net/articles/marie-estelle-dupont-psycho-couacs-a-droite-toute")
$wsArticles.Range("E10").Style = $wsArticles.Range("E9").Style

# --- Debunks sheet: point the existing debunk row to the new article, and fix text ---
$wsDebunks.Range("B2").Value = 9
$wsDebunks.Range("C2").Value = "Précision"

# --- Selections / active sheet bookkeeping ---
[void]$wsArticles.Range("F18").Select()

[void]$wsDebunks.Activate()
[void]$wsDebunks.Range("E12").Select()
